$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update text content
$ws.Range("P1").Value = "orientáció"
$ws.Range("G2").Value = "tudás és kihívások keresése"

# Center-align the two cells that now carry the shorter labels
$ws.Range("G2").HorizontalAlignment = -4108
$ws.Range("K2").HorizontalAlignment = -4108

# Adjust column widths (D, E, F, G)
$ws.Columns.Item(4).ColumnWidth = 6.42578125
$ws.Columns.Item(5).ColumnWidth = 8
$ws.Columns.Item(6).ColumnWidth = 20.28515625
$ws.Columns.Item(7).ColumnWidth = 25.85546875

# Update the last active selection
$ws.Range("M13").Select()
